# Femacal de La Calera - Haba: weekly price update.
# A new observation is inserted as row 36 (pushing the existing rows
# 36-73 down to 37-74), matching the author's "Fruta / hortaliza, semanal"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36, shifting rows 36:73 down to 37:74.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with this week's data point.
$ws.Range("A36").Value = 3
$ws.Range("B36").Value = "Femacal de La Calera"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 44484
$ws.Range("E36").Value = 5
$ws.Range("F36").Value = 100112026
$ws.Range("G36").Value = "Haba"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 105
$ws.Range("K36").Value = 8000
$ws.Range("L36").Value = 8500
$ws.Range("M36").Value = 8238
$ws.Range("N36").Value = "$/malla 25 kilos"
$ws.Range("O36").Value = "Provincia de Quillota"
$ws.Range("P36").Value = 330
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
